# Sync workbook with current readxl test sheets:
#  - drop the now-unused Sheet2 and Sheet3
#  - Sheet1's data actually lives one column to the left (A:D instead of B:E),
#    so remove the empty leading column A and let everything shift over.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A:A").EntireColumn.Delete()
